$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rotate the blog article rota forward: drop article 84, shift 85->I7, 86->E7,
# and bring article 87 live in C7.
$ws.Range("I7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 85"
$ws.Range("E7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 86"
$ws.Range("C7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 87"

# Scroll the sheet view back so column A (instead of B) is the left-most
# visible column, keeping row 6 as the top-most visible row.
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 6
